$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status")

# Wipe the old layout (values + formatting) so moved rows don't leave
# stale bold/blank cells behind in the gaps they vacate.
$ws.Range("B1:B20").Clear()

# Rows whose text is unchanged from the original sheet - write these
# first so any brand-new shared-string entries created below are
# appended after them, same as the source edit.
$unchanged = @(
    @{ Row = 2;  Bold = $true;  Text = "ITEMS TO COMPLETE" },
    @{ Row = 5;  Bold = $true;  Text = "function songPlay" },
    @{ Row = 6;  Bold = $false; Text = "needs to play one of the 4 songs by populating the audio tag in the html page" },
    @{ Row = 7;  Bold = $false; Text = "needs to randomly choose the song from the album array" },
    @{ Row = 15; Bold = $true;  Text = "function displayCorrectWrong" },
    @{ Row = 16; Bold = $false; Text = "displays a document.write depending if answer is correct or not" },
    @{ Row = 19; Bold = $true;  Text = "function logPoint" },
    @{ Row = 20; Bold = $false; Text = "will log a point if the correct album is clicked" },
    @{ Row = 31; Bold = $true;  Text = "function displayScores" },
    @{ Row = 32; Bold = $false; Text = "make this a function of the getWinner function" },
    @{ Row = 33; Bold = $false; Text = "how to pull a Var into a document.write function" }
)

foreach ($r in $unchanged) {
    $cell = $ws.Cells.Item($r.Row, 2)
    $cell.Value = $r.Text
    $cell.Font.Bold = $r.Bold
}

# New notes documenting the 6th-random-album no-dupe fix (albumClick /
# gameFinish / nextRound). Written in shared-string creation order:
# albumClick, the player function..., function gameFinish, function nextRound.
$new = @(
    @{ Row = 10; Bold = $true;  Text = "function albumClick" },
    @{ Row = 11; Bold = $false; Text = "the player function that selects the album" },
    @{ Row = 28; Bold = $true;  Text = "function gameFinish" },
    @{ Row = 23; Bold = $true;  Text = "function nextRound" }
)

foreach ($r in $new) {
    $cell = $ws.Cells.Item($r.Row, 2)
    $cell.Value = $r.Text
    $cell.Font.Bold = $r.Bold
}

$ws.Range("A21").Select()
